# Fixes during Regression Testing
#
# vRelayPaymentsACH.xlsx is a Katalon data-bootstrap workbook: each worksheet's
# row 2 (row 3 on "CMCAutopayPS", which carries two scenario rows) stores the
# outcome of the most recent regression run in column A ("Result": Pass/Fail)
# and column B ("Date": execution timestamp). Re-running the regression suite
# just overwrites those two cells per sheet with the latest run's outcome.
#
# This replays that update for every sheet whose last-run Result/Date changed.

$wb = $excel.ActiveWorkbook

function Set-RunOutcome {
    param(
        [string]$SheetName,
        [int]$Row,
        [string]$Result,
        [string]$Timestamp
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range("A$Row").Value = $Result
    $ws.Range("B$Row").Value = $Timestamp
}

Set-RunOutcome "PayNowNoCFPC"              2 "Pass" "Tue Nov 18 02:07:24 IST 2025"
Set-RunOutcome "PayNowNoCFPS"              2 "Pass" "Tue Nov 18 02:11:00 IST 2025"
Set-RunOutcome "PayNowNoCFCorp"            2 "Pass" "Tue Nov 18 02:06:28 IST 2025"
Set-RunOutcome "PayNowSCFPC"               2 "Fail" "Tue Nov 18 02:34:53 IST 2025"
Set-RunOutcome "PayNowSCFPS"               2 "Fail" "Tue Nov 18 02:42:57 IST 2025"
Set-RunOutcome "PayNowSCFCorp"             2 "Fail" "Tue Nov 18 02:26:52 IST 2025"
Set-RunOutcome "PayNowDCFPC"               2 "Fail" "Tue Nov 18 01:50:28 IST 2025"
Set-RunOutcome "PayNowDCFPS"               2 "Fail" "Tue Nov 18 01:58:50 IST 2025"
Set-RunOutcome "PayNowDCFCorp"             2 "Fail" "Tue Nov 18 01:41:49 IST 2025"
Set-RunOutcome "SCFPSVerbiage"             2 "Pass" "Tue Nov 18 02:51:15 IST 2025"
Set-RunOutcome "SCFPCVerbiage"             2 "Pass" "Tue Nov 18 02:50:34 IST 2025"
Set-RunOutcome "SCFCorpVerbiage"           2 "Pass" "Tue Nov 18 02:49:28 IST 2025"
Set-RunOutcome "DCFPSVerbiage"             2 "Pass" "Tue Nov 18 02:47:40 IST 2025"
Set-RunOutcome "DCFPCVerbiage"             2 "Pass" "Tue Nov 18 02:47:09 IST 2025"
Set-RunOutcome "DCFCorpVerbiage"           2 "Pass" "Tue Nov 18 02:46:26 IST 2025"
Set-RunOutcome "CMCAutopayPS"              3 "Fail" "Tue Nov 18 00:49:28 IST 2025"
Set-RunOutcome "CMCAutoPayPC"              2 "Pass" "Tue Nov 18 00:32:52 IST 2025"
Set-RunOutcome "CMCAutoPayCorp"            2 "Fail" "Tue Nov 18 00:30:59 IST 2025"
Set-RunOutcome "CCDeferredPS"              2 "Pass" "Tue Nov 18 01:15:56 IST 2025"
Set-RunOutcome "CCDeferredPC"              2 "Pass" "Tue Nov 18 01:14:37 IST 2025"
Set-RunOutcome "CCDeferredCorp"            2 "Pass" "Tue Nov 18 01:13:04 IST 2025"
Set-RunOutcome "NoModifyAmountPC"          2 "Pass" "Tue Nov 18 01:19:42 IST 2025"
Set-RunOutcome "NoModifyAmountPS"          2 "Pass" "Tue Nov 18 01:20:52 IST 2025"
Set-RunOutcome "NoModifyAmountCorp"        2 "Pass" "Tue Nov 18 01:23:14 IST 2025"
Set-RunOutcome "NoModifyBillingAddressPC"  2 "Pass" "Tue Nov 18 01:24:21 IST 2025"
Set-RunOutcome "NoModifyBillingAddressPS"  2 "Pass" "Tue Nov 18 01:25:29 IST 2025"
